# MinAvgTwoSlice worksheet: insert a new working column before F (this shifts
# the existing F/G/I/J/L columns one place to the right -> G/H/J/K/M in both
# tables), then populate the freshly inserted E/F columns with a new
# "difference" / "running total via difference" helper calculation, and add
# a few scratch literal values off to the right (O2:Q4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at F; everything from F rightwards (F,G,I,J,L ...)
# shifts one column to the right (-> G,H,J,K,M ...) in every row, including
# the second mini-table in rows 11-14.
$ws.Columns("F").Insert()

# New column E: row-over-row difference of the input column D.
$ws.Range("E3").Formula = "=D3-D2"
$ws.Range("E4:E8").Formula = "=D4-D3"

# New column F: running total built from the new diff column instead of D
# directly (F2 - the old AVERAGE formula that is now in G2 - stays the seed).
$ws.Range("F3").Formula = "=E3+F2"
$ws.Range("F4:F8").Formula = "=E4+F3"

# A few scratch literal values added near the first table.
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1
$ws.Range("P3").Value = 15
$ws.Range("Q3").Value = 14
$ws.Range("P4").Value = 5
$ws.Range("Q4").Value = -10

# Move the active selection to match where the author ended up.
$ws.Range("F8").Select()
